$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column F ("Precision@5"), shifting existing
# F:N columns (Precision@5 .. ILD@5) one column to the right.
$ws.Columns("F:F").Insert()

# New header for the inserted column
$ws.Range("F1").Value = "HitRate@5"

# New HitRate@5 values for each data row
$ws.Range("F2").Value = 0.165424739195231
$ws.Range("F3").Value = 0.09239940387481371
$ws.Range("F4").Value = 0.9985096870342772
